$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Internal Assignment" column (O) to the experiment-type table.
# Header cell (O4): same look as the other bold headers (K4/L4..N4) but one point larger.
$ws.Range("N4").Copy($ws.Range("O4"))
$ws.Range("O4").Value = "Internal Assignment"
$ws.Range("O4").Font.Size = 12

# Data cells (O5:O8): same TRUE/FALSE style as column L, value "FALSE".
$ws.Range("L5").Copy($ws.Range("O5"))
$ws.Range("L6").Copy($ws.Range("O6"))
$ws.Range("L7").Copy($ws.Range("O7"))
$ws.Range("L8").Copy($ws.Range("O8"))

$ws.Range("P9").Select()
